$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write string-valued cells in a specific first-seen order so the
# rebuilt sharedStrings table comes out in the same order the source diff expects ---

# header row (column labels) - establishes shared-string indices 0-7
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("A2").Value = "name"

# negative-word column, top to bottom (indices 8-36)
$ws.Range("A3").Value = "poorly"
$ws.Range("A4").Value = "disappointing"
$ws.Range("A5").Value = "disappointed"
$ws.Range("A6").Value = "however"
$ws.Range("A7").Value = "broke"
$ws.Range("A8").Value = "poor"
$ws.Range("A9").Value = "waste"
$ws.Range("A10").Value = "instead"
$ws.Range("A11").Value = "junk"
$ws.Range("A12").Value = "smaller"
$ws.Range("A13").Value = "small"
$ws.Range("A14").Value = "paint"
$ws.Range("A15").Value = "broken"
$ws.Range("A16").Value = "plastic"
$ws.Range("A17").Value = "apart"
$ws.Range("A18").Value = "ok"
$ws.Range("A19").Value = "thought"
$ws.Range("A20").Value = "cheap"
$ws.Range("A21").Value = "though"
$ws.Range("A22").Value = "money"
$ws.Range("A23").Value = "size"
$ws.Range("A24").Value = "item"
$ws.Range("A25").Value = "hard"
$ws.Range("A26").Value = "work"
$ws.Range("A27").Value = "would"
$ws.Range("A28").Value = "product"
$ws.Range("A29").Value = "price"
$ws.Range("A30").Value = "use"
$ws.Range("A31").Value = "like"

# negative-table title (index 37)
$ws.Range("A1").Value = "negative"

# positive-word column, top to bottom (indices 38-49)
$ws.Range("J3").Value = "wonderful"
$ws.Range("J4").Value = "awesome"
$ws.Range("J5").Value = "favorite"
$ws.Range("J6").Value = "classic"
$ws.Range("J7").Value = "excellent"
$ws.Range("J8").Value = "great"
$ws.Range("J9").Value = "love"
$ws.Range("J10").Value = "loves"
$ws.Range("J11").Value = "perfect"
$ws.Range("J12").Value = "loved"
$ws.Range("J13").Value = "fun"
$ws.Range("J14").Value = "game"

# positive-table title (index 50)
$ws.Range("J1").Value = "positive"

# remaining header cells (K2..Q2) reuse already-registered shared strings
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"

# --- Step 2: numeric / boolean data cells ---

# row 3
$ws.Range("B3").Value = 0.9565217391304348
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = 44
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 2
$ws.Range("K3").Value = 0.875
$ws.Range("L3").Value = 49
$ws.Range("M3").Value = 49
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 7

# row 4
$ws.Range("B4").Value = 0.8181818181818182
$ws.Range("C4").Value = 36
$ws.Range("D4").Value = 36
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 8
$ws.Range("K4").Value = 0.8307692307692308
$ws.Range("L4").Value = 54
$ws.Range("M4").Value = 54
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 11

# row 5
$ws.Range("B5").Value = 0.7473118279569892
$ws.Range("C5").Value = 139
$ws.Range("D5").Value = 139
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 47
$ws.Range("K5").Value = 0.6344086021505376
$ws.Range("L5").Value = 59
$ws.Range("M5").Value = 59
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 34

# row 6
$ws.Range("B6").Value = 0.734375
$ws.Range("C6").Value = 47
$ws.Range("D6").Value = 47
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 17
$ws.Range("K6").Value = 0.5849056603773585
$ws.Range("L6").Value = 31
$ws.Range("M6").Value = 31
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 22

# row 7
$ws.Range("B7").Value = 0.7233009708737864
$ws.Range("C7").Value = 149
$ws.Range("D7").Value = 149
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 57
$ws.Range("K7").Value = 0.515625
$ws.Range("L7").Value = 33
$ws.Range("M7").Value = 33
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 31

# row 8
$ws.Range("B8").Value = 0.6901408450704225
$ws.Range("C8").Value = 49
$ws.Range("D8").Value = 49
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 22
$ws.Range("K8").Value = 0.3434426229508197
$ws.Range("L8").Value = 419
$ws.Range("M8").Value = 419
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 801

# row 9
$ws.Range("B9").Value = 0.6486486486486487
$ws.Range("C9").Value = 96
$ws.Range("D9").Value = 96
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 52
$ws.Range("K9").Value = 0.3127690100430416
$ws.Range("L9").Value = 218
$ws.Range("M9").Value = 218
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 479

# row 10
$ws.Range("B10").Value = 0.625
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 18
$ws.Range("K10").Value = 0.2697095435684647
$ws.Range("L10").Value = 130
$ws.Range("M10").Value = 130
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 352

# row 11
$ws.Range("B11").Value = 0.5636363636363636
$ws.Range("C11").Value = 31
$ws.Range("D11").Value = 31
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 24
$ws.Range("K11").Value = 0.1867469879518072
$ws.Range("L11").Value = 31
$ws.Range("M11").Value = 31
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 135

# row 12
$ws.Range("B12").Value = 0.5462184873949579
$ws.Range("C12").Value = 65
$ws.Range("D12").Value = 65
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 54
$ws.Range("K12").Value = 0.1651376146788991
$ws.Range("L12").Value = 54
$ws.Range("M12").Value = 54
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 273

# row 13
$ws.Range("B13").Value = 0.4840579710144928
$ws.Range("C13").Value = 167
$ws.Range("D13").Value = 167
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 178
$ws.Range("K13").Value = 0.07099035933391762
$ws.Range("L13").Value = 81
$ws.Range("M13").Value = 81
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 1060

# row 14
$ws.Range("B14").Value = 0.4761904761904762
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 33
$ws.Range("K14").Value = 0.02401038286826736
$ws.Range("L14").Value = 37
$ws.Range("M14").Value = 37
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 1504

# row 15
$ws.Range("B15").Value = 0.4337349397590362
$ws.Range("C15").Value = 36
$ws.Range("D15").Value = 36
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 47

# row 16
$ws.Range("B16").Value = 0.4330708661417323
$ws.Range("C16").Value = 55
$ws.Range("D16").Value = 55
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 72

# row 17
$ws.Range("B17").Value = 0.3789473684210526
$ws.Range("C17").Value = 36
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 59

# row 18
$ws.Range("B18").Value = 0.328125
$ws.Range("C18").Value = 42
$ws.Range("D18").Value = 42
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 86

# row 19
$ws.Range("B19").Value = 0.3168316831683168
$ws.Range("C19").Value = 64
$ws.Range("D19").Value = 64
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 138

# row 20
$ws.Range("B20").Value = 0.2796208530805687
$ws.Range("C20").Value = 59
$ws.Range("D20").Value = 59
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 152

# row 21
$ws.Range("B21").Value = 0.2564102564102564
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 87

# row 22
$ws.Range("B22").Value = 0.1962025316455696
$ws.Range("C22").Value = 62
$ws.Range("D22").Value = 62
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 254

# row 23
$ws.Range("B23").Value = 0.1958762886597938
$ws.Range("C23").Value = 38
$ws.Range("D23").Value = 38
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 156

# row 24
$ws.Range("B24").Value = 0.1884057971014493
$ws.Range("C24").Value = 52
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 224

# row 25
$ws.Range("B25").Value = 0.185
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 37
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 163

# row 26
$ws.Range("B26").Value = 0.1740506329113924
$ws.Range("C26").Value = 55
$ws.Range("D26").Value = 55
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 261

# row 27
$ws.Range("B27").Value = 0.1661721068249258
$ws.Range("C27").Value = 112
$ws.Range("D27").Value = 112
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 562

# row 28
$ws.Range("B28").Value = 0.1409691629955947
$ws.Range("C28").Value = 64
$ws.Range("D28").Value = 64
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 390

# row 29
$ws.Range("B29").Value = 0.1239193083573487
$ws.Range("C29").Value = 43
$ws.Range("D29").Value = 44
$ws.Range("E29").Value = 0.02
$ws.Range("F29").Value = 0.98
$ws.Range("G29").Value = $true
$ws.Range("H29").Value = 304

# row 30
$ws.Range("B30").Value = 0.08493150684931507
$ws.Range("C30").Value = 31
$ws.Range("D30").Value = 31
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 334

# row 31
$ws.Range("B31").Value = 0.05436573311367381
$ws.Range("C31").Value = 33
$ws.Range("D31").Value = 34
$ws.Range("E31").Value = 0.03
$ws.Range("F31").Value = 0.97
$ws.Range("G31").Value = $true
$ws.Range("H31").Value = 574

# --- Step 3: apply the bold/centered/bordered header style to the two new rows ---
# (copies the format already used by A3:A29 so the style table is reused, not duplicated)
$ws.Range("A30").Font.Bold = $true
$ws.Range("A30").HorizontalAlignment = -4108
$ws.Range("A30").VerticalAlignment = -4160
$ws.Range("A30").Borders.LineStyle = 1
$ws.Range("A31").Font.Bold = $true
$ws.Range("A31").HorizontalAlignment = -4108
$ws.Range("A31").VerticalAlignment = -4160
$ws.Range("A31").Borders.LineStyle = 1
